$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people) - row 13
$ws.Range("B13").Value = "'36.98"
$ws.Range("C13").Value = "'2.22"
$ws.Range("D13").Value = "'39.21"

# Employment (% of total) - row 14
$ws.Range("B14").Value = "'28.65"
$ws.Range("C14").Value = "'33.96"
$ws.Range("D14").Value = "'62.61"

# Enterprises (% of total) - row 16
$ws.Range("B16").Value = "'94.17"
$ws.Range("C16").Value = "'5.66"
$ws.Range("D16").Value = "'99.83"
